$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '''68.530.37'
$ws.Range('D2').Style = "Normal"
$ws.Range('E2').Value = '  +1.82%  '
$ws.Range('D3').Value = '''3.570.04'
$ws.Range('D3').Style = "Normal"
$ws.Range('E3').Value = '  +2.19%  '
$ws.Range('E4').Value = '  +0.02%  '
$ws.Range('D5').Value = '''622.57'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  +2.85%  '
$ws.Range('D6').Value = '''155.29'
$ws.Range('D6').Style = "Normal"
$ws.Range('E6').Value = '  +4.50%  '
$ws.Range('D7').Value = '''3.568.55'
$ws.Range('D7').Style = "Normal"
$ws.Range('E7').Value = '  +2.16%  '
$ws.Range('E8').Value = '  -0.04%  '
$ws.Range('D9').Value = '''0.492'
$ws.Range('D9').Style = "Normal"
$ws.Range('E9').Value = '  +2.31%  '
$ws.Range('E10').Value = '  +5.72%  '
$ws.Range('D11').Value = '''7.36'
$ws.Range('D11').Style = "Normal"
$ws.Range('E11').Value = '  +5.67%  '
$ws.Range('D12').Value = '''0.439'
$ws.Range('D12').Style = "Normal"
$ws.Range('E12').Value = '  +3.99%  '
$ws.Range('D13').Value = '''0.0000222'
$ws.Range('D13').Style = "Normal"
$ws.Range('E13').Value = '  +1.67%  '
$ws.Range('D14').Value = '''33.21'
$ws.Range('D14').Style = "Normal"
$ws.Range('E14').Value = '  +5.71%  '
$ws.Range('D15').Value = '''4.172.71'
$ws.Range('D15').Style = "Normal"
$ws.Range('E15').Value = '  +2.17%  '
$ws.Range('D16').Value = '''3.569.90'
$ws.Range('D16').Style = "Normal"
$ws.Range('E16').Value = '  +2.01%  '
$ws.Range('D17').Value = '''68.620.18'
$ws.Range('D17').Style = "Normal"
$ws.Range('E17').Value = '  +2.13%  '
$ws.Range('E18').Value = '  +0.01%  '
$ws.Range('D19').Value = '''6.78'
$ws.Range('D19').Style = "Normal"
$ws.Range('E19').Value = '  +6.13%  '
$ws.Range('D20').Value = '''16.01'
$ws.Range('D20').Style = "Normal"
$ws.Range('E20').Value = '  +6.67%  '
$ws.Range('D21').Value = '''9.99'
$ws.Range('D21').Style = "Normal"
$ws.Range('E21').Value = '  +10.81%  '
$ws.Range('D22').Value = '''456.71'
$ws.Range('D22').Style = "Normal"
$ws.Range('E22').Value = '  +2.38%  '
$ws.Range('E23').Value = '  +3.73%  '
$ws.Range('D24').Value = '''78.78'
$ws.Range('D24').Style = "Normal"
$ws.Range('E24').Value = '  +2.13%  '
$ws.Range('E25').Value = '  +2.41%  '
$ws.Range('D26').Value = '''3.709.06'
$ws.Range('D26').Style = "Normal"
$ws.Range('E26').Value = '  +2.11%  '
$ws.Range('E27').Value = '  -0.09%  '
$ws.Range('D28').Value = '''10.57'
$ws.Range('D28').Style = "Normal"
$ws.Range('E28').Value = '  +4.27%  '
$ws.Range('D29').Value = '''9.14'
$ws.Range('D29').Style = "Normal"
$ws.Range('E29').Value = '  +10.54%  '
$ws.Range('D30').Value = '''1.71'
$ws.Range('D30').Style = "Normal"
$ws.Range('E30').Value = '  +9.61%  '
$ws.Range('E31').Value = '  +3.94%  '
$ws.Range('D32').Value = '''0.171'
$ws.Range('D32').Style = "Normal"
$ws.Range('E32').Value = '  +5.15%  '
$ws.Range('E33').Value = '  -0.10%  '
$ws.Range('B34').Value = 'EthereumClassic'
$ws.Range('C34').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D34').Value = '''26.18'
$ws.Range('D34').Style = "Normal"
$ws.Range('E34').Value = '  +2.10%  '
$ws.Range('B35').Value = 'NEARProtocol'
$ws.Range('C35').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D35').Value = '''6.35'
$ws.Range('D35').Style = "Normal"
$ws.Range('E35').Value = '  +3.55%  '
$ws.Range('E36').Value = '  +3.75%  '
$ws.Range('D37').Value = '''3.563.09'
$ws.Range('D37').Style = "Normal"
$ws.Range('E37').Value = '  +2.31%  '
$ws.Range('E38').Value = '  +3.76%  '
$ws.Range('D39').Value = '''2.39'
$ws.Range('D39').Style = "Normal"
$ws.Range('E39').Value = '  +9.57%  '
$ws.Range('E40').Value = '  -0.03%  '
$ws.Range('D41').Value = '''178.95'
$ws.Range('D41').Style = "Normal"
$ws.Range('E41').Value = '  +3.06%  '
$ws.Range('D42').Value = '''0.0921'
$ws.Range('D42').Style = "Normal"
$ws.Range('E42').Value = '  +5.18%  '
$ws.Range('E43').Value = '  +0.01%  '
$ws.Range('D44').Value = '''5.61'
$ws.Range('D44').Style = "Normal"
$ws.Range('E44').Value = '  +3.61%  '
$ws.Range('D45').Value = '''31.07'
$ws.Range('D45').Style = "Normal"
$ws.Range('E45').Value = '  +15.26%  '
$ws.Range('E46').Value = '  +2.00%  '
$ws.Range('D47').Value = '''46.52'
$ws.Range('D47').Style = "Normal"
$ws.Range('E47').Value = '  +2.36%  '
$ws.Range('E48').Value = '  +6.82%  '
$ws.Range('D49').Value = '''2.68'
$ws.Range('D49').Style = "Normal"
$ws.Range('E49').Value = '  +4.09%  '
$ws.Range('D50').Value = '''7.81'
$ws.Range('D50').Style = "Normal"
$ws.Range('E50').Value = '  +3.80%  '
$ws.Range('D51').Value = '''0.264'
$ws.Range('D51').Style = "Normal"
$ws.Range('E51').Value = '  +7.92%  '
